$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Sheet1 (Recommandations) data: row -> (A,B,C,D,E,F,G)
$sheet1Data = @(
    @(2, "UNIWAX CI", 0, 4, 3590, 930, "🟡 Observer", "➖ Neutre"),
    @(3, "CFAO MOTORS CI", 0, 4, 3530, 975, "🟡 Observer", "➖ Neutre"),
    @(4, "BRVM - SERVICES PUBLICS", 0, 8, 3312.86, 104, "🟡 Observer", "➖ Neutre"),
    @(5, "SETAO CI", 0, 4, 2780, 710, "🟡 Observer", "➖ Neutre"),
    @(6, "NEI-CEDA CI", 0, 4, 2560, 655, "🟡 Observer", "➖ Neutre"),
    @(7, "BRVM - AUTRES SECTEURS", 0, 4, 2512.27, 620.37, "🟡 Observer", "➖ Neutre"),
    @(8, "AIR LIQUIDE CI", 0, 4, 2390, 595, "🟡 Observer", "➖ Neutre"),
    @(9, "BRVM - DISTRIBUTION", 0, 4, 1584.19, 404.87, "🟡 Observer", "➖ Neutre"),
    @(10, "BRVM - TRANSPORT", 0, 4, 1382.96, 342.68, "🟡 Observer", "➖ Neutre"),
    @(11, "BRVM - AGRICULTURE", 0, 4, 1319.14, 329.26, "🟡 Observer", "➖ Neutre"),
    @(12, "BRVM - INDUSTRIELS", 0, 4, 550.37, 136.53, "🟡 Observer", "➖ Neutre"),
    @(13, "BRVM-PRESTIGE", 0, 4, 532.38, 133.43, "🟡 Observer", "➖ Neutre"),
    @(14, "BRVM - FINANCES", 0, 4, 502.99, 126.28, "🟡 Observer", "➖ Neutre"),
    @(15, "BRVM - SERVICES FINANCIERS", 0, 4, 494.32, 124.1, "🟡 Observer", "➖ Neutre"),
    @(16, "BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 4, 489.06, 127.15, "🟡 Observer", "➖ Neutre"),
    @(17, "BRVM - ENERGIE", 0, 4, 435.33, 106.89, "🟡 Observer", "➖ Neutre"),
    @(18, "BRVM - CONSOMMATION DE BASE          (**)", 0, 2, 419.67, 210.16, "🟡 Observer", "➖ Neutre"),
    @(19, "BRVM - TELECOMMUNICATIONS", 0, 4, 378.19, 94.70999999999999, "🟡 Observer", "➖ Neutre"),
    @(20, "BRVM - INDUSTRIE       (**)", 0, 1, 262.79, 262.79, "🟡 Observer", "➖ Neutre"),
    @(21, "BRVM - INDUSTRIE                  (**)", 0, 1, 251.95, 251.95, "🟡 Observer", "➖ Neutre"),
    @(22, "BRVM - INDUSTRIE             (**)", 0, 1, 251.86, 251.86, "🟡 Observer", "➖ Neutre"),
    @(23, "BRVM - CONSOMMATION DE BASE      (**)", 0, 1, 217.36, 217.36, "🟡 Observer", "➖ Neutre"),
    @(24, "BRVM-PRINCIPAL          (**)", 0, 1, 192.79, 192.79, "🟡 Observer", "➖ Neutre"),
    @(25, "BRVM-PRINCIPAL                (**)", 0, 1, 191.84, 191.84, "🟡 Observer", "➖ Neutre"),
    @(26, "BRVM-PRINCIPAL                    (**)", 0, 1, 191.6, 191.6, "🟡 Observer", "➖ Neutre"),
    @(27, "CFAO MOTORS CI (CFAC)", 4, 0, 27.45, 6.04, "🟢 Achat", "✅ Renforcer"),
    @(28, "ORAGROUP TOGO (ORGT)", 2, 0, 8.9, 5.57, "🟡 Observer", "➖ Neutre"),
    @(29, "NEI-CEDA CI (NEIC)", 1, 0, 6.98, 6.98, "🟡 Observer", "➖ Neutre"),
    @(30, "UNIWAX CI (UNXC)", 2, 1, 5.48, -6.45, "🟡 Observer", "👀 À surveiller"),
    @(31, "TOTALENERGIES MARKETING CI (TTLC)", 1, 0, 4.26, 4.26, "🟡 Observer", "➖ Neutre"),
    @(32, "BANK OF AFRICA BN (BOAB)", 1, 0, 4.01, 4.01, "🟡 Observer", "➖ Neutre"),
    @(33, "NESTLE CI (NTLC)", 1, 0, 3.34, 3.34, "🟡 Observer", "➖ Neutre"),
    @(34, "TOTALENERGIES MARKETING SN (TTLS)", 1, 0, 3.2, 3.2, "🟡 Observer", "➖ Neutre"),
    @(35, "SETAO CI (STAC)", 1, 2, 2.5, -2.11, "🟡 Observer", "👀 À surveiller"),
    @(36, "BANK OF AFRICA ML (BOAM)", 1, 0, 2.23, 2.23, "🟡 Observer", "➖ Neutre"),
    @(37, "SOCIETE GENERALE COTE D'IVOIRE (SGBC)", 1, 0, 1.96, 1.96, "🟡 Observer", "➖ Neutre"),
    @(38, "ECOBANK TRANS. INCORP. TG (ETIT)", 1, 1, 0.37, 6.25, "🟡 Observer", "👀 À surveiller"),
    @(39, "ONATEL BF (ONTBF)", 1, 1, 0.27, 2.78, "🟡 Observer", "👀 À surveiller"),
    @(40, "TOTAL", 0, 3, 0, 0, "🟡 Observer", "➖ Neutre"),
    @(41, "BANK OF AFRICA NG (BOAN)", 1, 1, -0.9399999999999999, 4, "🟡 Observer", "👀 À surveiller"),
    @(42, "BANK OF AFRICA BF (BOABF)", 0, 1, -1.05, -1.05, "🟡 Observer", "➖ Neutre"),
    @(43, "AFRICA GLOBAL LOGISTICS CI (SDSC)", 0, 1, -1.06, -1.06, "🟡 Observer", "➖ Neutre"),
    @(44, "SICABLE CI (CABC)", 0, 1, -1.13, -1.13, "🟡 Observer", "➖ Neutre"),
    @(45, "FILTISAC CI (FTSC)", 0, 1, -1.24, -1.24, "🟡 Observer", "➖ Neutre"),
    @(46, "VIVO ENERGY CI (SHEC)", 0, 1, -1.3, -1.3, "🟡 Observer", "➖ Neutre"),
    @(47, "SOGB CI (SOGC)", 0, 1, -2.47, -2.47, "🟡 Observer", "➖ Neutre"),
    @(48, "TRACTAFRIC MOTORS CI (PRSC)", 0, 1, -3.05, -3.05, "🟡 Observer", "➖ Neutre"),
    @(49, "BERNABE CI (BNBC)", 0, 2, -3.07, -1.96, "🟡 Observer", "➖ Neutre"),
    @(50, "SERVAIR ABIDJAN CI (ABJC)", 0, 1, -3.31, -3.31, "🟡 Observer", "➖ Neutre"),
    @(51, "SMB CI (SMBC)", 0, 1, -6.88, -6.88, "🟡 Observer", "➖ Neutre"),
    @(52, "SICOR CI (SICC)", 0, 1, -7.38, -7.38, "🟡 Observer", "➖ Neutre"),
    @(53, "UNILEVER CI (UNLC)", 0, 2, -14.46, -7.49, "🟡 Observer", "➖ Neutre"),
)

# Sheet2 (Top_YTD) data: row -> (A,B)
$sheet2Data = @(
    @(2, "BRVM - SERVICES PUBLICS", 7949652.55),
    @(3, "UNIWAX CI", 986752.53),
    @(4, "CFAO MOTORS CI", 923058.9399999999),
    @(5, "SETAO CI", 398644.78),
    @(6, "NEI-CEDA CI", 299628.87),
    @(7, "BRVM - AUTRES SECTEURS", 280865.61),
    @(8, "AIR LIQUIDE CI", 236582.25),
    @(9, "BRVM - DISTRIBUTION", 60428.07),
    @(10, "BRVM - TRANSPORT", 39371.67),
    @(11, "BRVM - AGRICULTURE", 34019.54),
)
foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
    $ws1.Cells.Item($r, 7).Value = $row[7]
}

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
}

Write-Output "Update complete. Sheet1 dimension: $($ws1.UsedRange.Address())"
